$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("hotel_info")

# Fill in previously-empty review-count / rank cells on row 2.
# The source data stores these as text (not numbers), so build each literal
# via TEXT() in a scratch cell, then copy/paste-special-values it into place.
# (A direct string assignment like Range("G2").Value = "5" gets silently
# re-coerced to a number since it looks numeric - going through a genuine
# text-typed formula result and a values-only paste keeps it text without
# picking up an extra NumberFormat/quote-prefix cell style.)
$scratch = $ws.Range("ZZ1")

$scratch.Formula = '=TEXT(5,"0")'
$scratch.Copy()
$ws.Range("G2").PasteSpecial(-4163)

$scratch.Formula = '=TEXT(312,"0")'
$scratch.Copy()
$ws.Range("H2").PasteSpecial(-4163)

$scratch.Formula = '=TEXT(5,"0")'
$scratch.Copy()
$ws.Range("I2").PasteSpecial(-4163)

$scratch.ClearContents()
$excel.CutCopyMode = 0
